$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.139906333333333
$ws.Range("H2").Value = 3.419719
$ws.Range("I2").Value = 0.2178538649973528
$ws.Range("J2").Value = 0.2178538649973527
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 55.56529466666667
$ws.Range("N2").Value = 166.695884
$ws.Range("O2").Value = 0.424029640296873
$ws.Range("P2").Value = 0.4240296402968731
$ws.Range("Q2").Value = 63.33923130406622
$ws.Range("R2").Value = 570.053081736596
$ws.Range("S2").Value = 0.09237649601211102
$ws.Range("T2").Value = 0.09237649601211102
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.139906333333333
$ws.Range("H3").Value = 3.419719
$ws.Range("I3").Value = 0.2178538649973528
$ws.Range("J3").Value = 0.2178538649973527
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 10.84181733333333
$ws.Range("N3").Value = 32.525452
$ws.Range("O3").Value = 0.08273603031526086
$ws.Range("P3").Value = 0.08273603031526089
$ws.Range("Q3").Value = 12.35865624310978
$ws.Range("R3").Value = 111.227906187988
$ws.Range("S3").Value = 0.01802436397871773
$ws.Range("T3").Value = 0.01802436397871773
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.139906333333333
$ws.Range("H4").Value = 3.419719
$ws.Range("I4").Value = 0.2178538649973528
$ws.Range("J4").Value = 0.2178538649973527
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 40.13853466666666
$ws.Range("N4").Value = 120.415604
$ws.Range("O4").Value = 0.3063050150071534
$ws.Range("P4").Value = 0.3063050150071535
$ws.Range("Q4").Value = 45.75416987725288
$ws.Range("R4").Value = 411.787528895276
$ws.Range("S4").Value = 0.0667297313873805
$ws.Range("T4").Value = 0.0667297313873805
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.139906333333333
$ws.Range("H5").Value = 3.419719
$ws.Range("I5").Value = 0.2178538649973528
$ws.Range("J5").Value = 0.2178538649973527
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.88776266666666
$ws.Range("N5").Value = 65.66328799999999
$ws.Range("O5").Value = 0.167029801355803
$ws.Range("P5").Value = 0.167029801355803
$ws.Range("Q5").Value = 24.94999928623022
$ws.Range("R5").Value = 224.549993576072
$ws.Range("S5").Value = 0.03638808779510175
$ws.Range("T5").Value = 0.03638808779510175
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.139906333333333
$ws.Range("H6").Value = 3.419719
$ws.Range("I6").Value = 0.2178538649973528
$ws.Range("J6").Value = 0.2178538649973527
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.607653333333333
$ws.Range("N6").Value = 7.82296
$ws.Range("O6").Value = 0.01989951302490964
$ws.Range("P6").Value = 0.01989951302490964
$ws.Range("Q6").Value = 2.972480549804444
$ws.Range("R6").Value = 26.75232494824
$ws.Range("S6").Value = 0.004335185824041727
$ws.Range("T6").Value = 0.004335185824041726
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.902924
$ws.Range("H7").Value = 5.708772
$ws.Range("I7").Value = 0.3636784322304457
$ws.Range("J7").Value = 0.3636784322304456
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 55.56529466666667
$ws.Range("N7").Value = 166.695884
$ws.Range("O7").Value = 0.424029640296873
$ws.Range("P7").Value = 0.4240296402968731
$ws.Range("Q7").Value = 105.736532788272
$ws.Range("R7").Value = 951.628795094448
$ws.Range("S7").Value = 0.1542104348024066
$ws.Range("T7").Value = 0.1542104348024066
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.902924
$ws.Range("H8").Value = 5.708772
$ws.Range("I8").Value = 0.3636784322304457
$ws.Range("J8").Value = 0.3636784322304456
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 10.84181733333333
$ws.Range("N8").Value = 32.525452
$ws.Range("O8").Value = 0.08273603031526086
$ws.Range("P8").Value = 0.08273603031526089
$ws.Range("Q8").Value = 20.631154407216
$ws.Range("R8").Value = 185.680389664944
$ws.Range("S8").Value = 0.03008930979402469
$ws.Range("T8").Value = 0.0300893097940247
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.902924
$ws.Range("H9").Value = 5.708772
$ws.Range("I9").Value = 0.3636784322304457
$ws.Range("J9").Value = 0.3636784322304456
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.13853466666666
$ws.Range("N9").Value = 120.415604
$ws.Range("O9").Value = 0.3063050150071534
$ws.Range("P9").Value = 0.3063050150071535
$ws.Range("Q9").Value = 76.38058094203198
$ws.Range("R9").Value = 687.425228478288
$ws.Range("S9").Value = 0.1113965276421247
$ws.Range("T9").Value = 0.1113965276421247
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.902924
$ws.Range("H10").Value = 5.708772
$ws.Range("I10").Value = 0.3636784322304457
$ws.Range("J10").Value = 0.3636784322304456
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 21.88776266666666
$ws.Range("N10").Value = 65.66328799999999
$ws.Range("O10").Value = 0.167029801355803
$ws.Range("P10").Value = 0.167029801355803
$ws.Range("Q10").Value = 41.65074888470399
$ws.Range("R10").Value = 374.856739962336
$ws.Range("S10").Value = 0.06074513629284119
$ws.Range("T10").Value = 0.06074513629284119
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.902924
$ws.Range("H11").Value = 5.708772
$ws.Range("I11").Value = 0.3636784322304457
$ws.Range("J11").Value = 0.3636784322304456
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.607653333333333
$ws.Range("N11").Value = 7.82296
$ws.Range("O11").Value = 0.01989951302490964
$ws.Range("P11").Value = 0.01989951302490964
$ws.Range("Q11").Value = 4.962166111679999
$ws.Range("R11").Value = 44.65949500512
$ws.Range("S11").Value = 0.007237023699048471
$ws.Range("T11").Value = 0.007237023699048469
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.189605333333333
$ws.Range("H12").Value = 6.568816
$ws.Range("I12").Value = 0.4184677027722017
$ws.Range("J12").Value = 0.4184677027722016
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 55.56529466666667
$ws.Range("N12").Value = 166.695884
$ws.Range("O12").Value = 0.424029640296873
$ws.Range("P12").Value = 0.4240296402968731
$ws.Range("Q12").Value = 121.6660655503715
$ws.Range("R12").Value = 1094.994589953344
$ws.Range("S12").Value = 0.1774427094823554
$ws.Range("T12").Value = 0.1774427094823555
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.189605333333333
$ws.Range("H13").Value = 6.568816
$ws.Range("I13").Value = 0.4184677027722017
$ws.Range("J13").Value = 0.4184677027722016
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 10.84181733333333
$ws.Range("N13").Value = 32.525452
$ws.Range("O13").Value = 0.08273603031526086
$ws.Range("P13").Value = 0.08273603031526089
$ws.Range("Q13").Value = 23.73930105609244
$ws.Range("R13").Value = 213.653709504832
$ws.Range("S13").Value = 0.03462235654251845
$ws.Range("T13").Value = 0.03462235654251846
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.189605333333333
$ws.Range("H14").Value = 6.568816
$ws.Range("I14").Value = 0.4184677027722017
$ws.Range("J14").Value = 0.4184677027722016
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 40.13853466666666
$ws.Range("N14").Value = 120.415604
$ws.Range("O14").Value = 0.3063050150071534
$ws.Range("P14").Value = 0.3063050150071535
$ws.Range("Q14").Value = 87.88754957831821
$ws.Range("R14").Value = 790.987946204864
$ws.Range("S14").Value = 0.1281787559776482
$ws.Range("T14").Value = 0.1281787559776482
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.189605333333333
$ws.Range("H15").Value = 6.568816
$ws.Range("I15").Value = 0.4184677027722017
$ws.Range("J15").Value = 0.4184677027722016
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 21.88776266666666
$ws.Range("N15").Value = 65.66328799999999
$ws.Range("O15").Value = 0.167029801355803
$ws.Range("P15").Value = 0.167029801355803
$ws.Range("Q15").Value = 47.92556186966755
$ws.Range("R15").Value = 431.330056827008
$ws.Range("S15").Value = 0.06989657726786005
$ws.Range("T15").Value = 0.06989657726786005
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.189605333333333
$ws.Range("H16").Value = 6.568816
$ws.Range("I16").Value = 0.4184677027722017
$ws.Range("J16").Value = 0.4184677027722016
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2.607653333333333
$ws.Range("N16").Value = 7.82296
$ws.Range("O16").Value = 0.01989951302490964
$ws.Range("P16").Value = 0.01989951302490964
$ws.Range("Q16").Value = 5.709731646151111
$ws.Range("R16").Value = 51.38758481536
$ws.Range("S16").Value = 0.008327303501819443
$ws.Range("T16").Value = 0.008327303501819441
